$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 347, pushing the existing rows
# 347-404 down to 349-406 (this mirrors the diff, which is a pure
# two-row insertion in the middle of the data block).
$ws.Range("A347:A348").EntireRow.Insert()

# --- Row 347: new data row (date 45077 / Primera) ---
$ws.Range("A347").Value = 7
$ws.Range("B347").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C347").Value = "Ñuble"
$ws.Range("D347").Value = 45077
$ws.Range("E347").Value = 16
$ws.Range("F347").Value = 100112006
$ws.Range("G347").Value = "Repollo"
$ws.Range("H347").Value = "Crespo record"
$ws.Range("I347").Value = "Primera"
$ws.Range("J347").Value = 500
$ws.Range("K347").Value = 1200
$ws.Range("L347").Value = 1300
$ws.Range("M347").Value = 1250
$ws.Range("N347").Value = "$/unidad"
$ws.Range("O347").Value = "Provincia de Diguillín"
$ws.Range("P347").Value = 1250
$ws.Range("Q347").Value = 1
$ws.Range("R347").Value = "Hortaliza"

# --- Row 348: new data row (date 45077 / Segunda) ---
$ws.Range("A348").Value = 7
$ws.Range("B348").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C348").Value = "Ñuble"
$ws.Range("D348").Value = 45077
$ws.Range("E348").Value = 16
$ws.Range("F348").Value = 100112006
$ws.Range("G348").Value = "Repollo"
$ws.Range("H348").Value = "Crespo record"
$ws.Range("I348").Value = "Segunda"
$ws.Range("J348").Value = 300
$ws.Range("K348").Value = 1000
$ws.Range("L348").Value = 1000
$ws.Range("M348").Value = 1000
$ws.Range("N348").Value = "$/unidad"
$ws.Range("O348").Value = "Provincia de Diguillín"
$ws.Range("P348").Value = 1000
$ws.Range("Q348").Value = 1
$ws.Range("R348").Value = "Hortaliza"
